{"js": "// \"What is synchronous para 2\" \u2014 split the intro paragraph so that the\n// standalone word \"synchronous\" is wrapped in a bookmark (as Word does when\n// pasting / DDE-linking a selection), then append a new paragraph of body\n// text after it, surrounded by blank spacer paragraphs (matching the\n// existing blank-paragraph spacing already used in the doc).\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The paragraph that begins \"Before writing asynchronous ...\" is the third\n// paragraph in the document (index 2): title, blank spacer, then this one.\nconst introPara = paragraphs.items[2];\n\n// Find the standalone \"synchronous\" occurrence (not the one embedded inside\n// \"asynchronous\"). matchWholeWord excludes the \"...a[synchronous]\" hit, so\n// this returns the single real word in the document.\nconst matches = context.document.body.search(\"synchronous\", { matchCase: false, matchWholeWord: true });\nmatches.load(\"text\");\nawait context.sync();\n\nconst synchronousRange = matches.items[0];\nsynchronousRange.insertBookmark(\"__DdeLink__0_1853420875\");\n\n// Append the new content after the intro paragraph: blank paragraph, the\n// new explanatory paragraph, then two trailing blank paragraphs \u2014 mirroring\n// the paragraph/run formatting (jc=\"both\", sz/szCs=26, u=\"none\") already in\n// use for this paragraph, which insertParagraph naturally inherits.\nconst blank1 = introPara.insertParagraph(\"\", Word.InsertLocation.after);\nconst bodyPara = blank1.insertParagraph(\n  \"In context of computer execution of process or task perform on threads. \" +\n  \"A thread is a series of commands that exists as a unit of work. The \" +\n  \"operating system simulates this by allocating slices of time to \" +\n  \"different threads. Now as days we have multi core and multi threads \" +\n  \"CPU, now doing task/processes synchronous will waste resources of the \" +\n  \"machine, because we can do multiple tasks simultaneously.\",\n  Word.InsertLocation.after\n);\nconst blank2 = bodyPara.insertParagraph(\"\", Word.InsertLocation.after);\nblank2.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# \"What is synchronous para 2\" \u2014 wrap the standalone word \"synchronous\" in\n# the intro paragraph with a bookmark (as Word does for a DDE-linked\n# selection), then add a new explanatory paragraph after it, padded with\n# blank spacer paragraphs to match the document's existing spacing style.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: bookmark the standalone \"synchronous\" ------------------------\n# MatchWholeWord avoids matching the \"synchronous\" substring inside\n# \"asynchronous\" earlier in the same paragraph; searching the whole\n# document Content finds the single real occurrence.\n$target = $d.Content\n$found = $target.Find.Execute(\"synchronous\", $false, $true)\nif ($found) {\n    $d.Bookmarks.Add(\"__DdeLink__0_1853420875\", $target) | Out-Null\n}\n\n# --- Step 2: append new paragraphs after the intro paragraph --------------\n$bodyText = \"In context of computer execution of process or task perform on threads. A thread is a series of commands that exists as a unit of work. The operating system simulates this by allocating slices of time to different threads. Now as days we have multi core and multi threads CPU, now doing task/processes synchronous will waste resources of the machine, because we can do multiple tasks simultaneously.\"\n\n$introIndex = 3\n$introPara = $d.Paragraphs.Item($introIndex)\n\n# Blank paragraph right after the intro paragraph.\n$r = $introPara.Range.Duplicate\n$r.Collapse(0)              # wdCollapseEnd\n$r.InsertParagraphAfter()\n\n# Paragraph holding the new explanatory text.\n$blank1 = $d.Paragraphs.Item($introIndex + 1)\n$r = $blank1.Range.Duplicate\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n$bodyPara = $d.Paragraphs.Item($introIndex + 2)\n$r = $bodyPara.Range.Duplicate\n$r.Collapse(0)\n$r.InsertAfter($bodyText)\n\n# Two trailing blank paragraphs.\n$r = $bodyPara.Range.Duplicate\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n$blank2 = $d.Paragraphs.Item($introIndex + 3)\n$r = $blank2.Range.Duplicate\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n"}
